$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.972.92'
$ws.Range("E2").Value = '  -0.76%  '
$ws.Range("D3").Value = '2.462.37'
$ws.Range("E3").Value = '  -1.26%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '558.96'
$ws.Range("E5").Value = '  -1.76%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '163.06'
$ws.Range("E6").Value = '  -1.86%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.504'
$ws.Range("E8").Value = '  -1.38%  '
$ws.Range("D9").Value = '2.461.16'
$ws.Range("E9").Value = '  -1.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.152'
$ws.Range("E10").Value = '  -4.79%  '
$ws.Range("E11").Value = '  -0.56%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.335'
$ws.Range("E12").Value = '  -3.63%  '
$ws.Range("E13").Value = '  -1.45%  '
$ws.Range("D14").Value = '2.917.13'
$ws.Range("E14").Value = '  -1.12%  '
$ws.Range("D15").Value = '68.988.14'
$ws.Range("E15").Value = '  -0.82%  '
$ws.Range("E16").Value = '  -3.42%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '23.60'
$ws.Range("E17").Value = '  -2.31%  '
$ws.Range("D18").Value = '2.461.47'
$ws.Range("E18").Value = '  -3.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.79'
$ws.Range("E19").Value = '  -3.62%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '341.18'
$ws.Range("E20").Value = '  -3.32%  '
$ws.Range("E21").Value = '  -5.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.80'
$ws.Range("E22").Value = '  -2.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.93'
$ws.Range("E23").Value = '  +0.79%  '
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '67.13'
$ws.Range("E25").Value = '  -3.24%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.69'
$ws.Range("E26").Value = '  -2.82%  '
$ws.Range("D27").Value = '2.590.14'
$ws.Range("E27").Value = '  -1.27%  '
$ws.Range("E28").Value = '  +0.23%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.20'
$ws.Range("E29").Value = '  -5.21%  '
$ws.Range("D30").Value = '0.0₃0820'
$ws.Range("E30").Value = '  -5.78%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.16'
$ws.Range("E31").Value = '  -4.99%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '439.09'
$ws.Range("E32").Value = '  -0.09%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.15'
$ws.Range("E34").Value = '  -3.86%  '
$ws.Range("E35").Value = '  -5.56%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '157.07'
$ws.Range("E36").Value = '  +2.57%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.04'
$ws.Range("E37").Value = '  -0.13%  '
$ws.Range("E38").Value = '  -0.03%  '
$ws.Range("E39").Value = '  -4.66%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.85'
$ws.Range("E40").Value = '  -1.59%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.46'
$ws.Range("E41").Value = '  -2.76%  '
$ws.Range("E42").Value = '  -4.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '37.46'
$ws.Range("E43").Value = '  -1.21%  '
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.48'
$ws.Range("E44").Value = '  -5.97%  '
$ws.Range("B45").Value = 'ImmutableX'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.11'
$ws.Range("E45").Value = '  +3.26%  '
$ws.Range("E46").Value = '  -3.92%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '133.35'
$ws.Range("E47").Value = '  -4.20%  '
$ws.Range("E48").Value = '  -2.51%  '
$ws.Range("E49").Value = '  -0.64%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.482'
$ws.Range("E50").Value = '  -4.40%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.561'
$ws.Range("E51").Value = '  -2.25%  '
